# Test data from invalid login page
#
# Adds two new worksheets ("valid_login" and "invalid_login") after the
# existing "script1" sheet, populates them with username/password test
# data, selects the relevant cells, and makes "invalid_login" the active
# (selected) sheet/tab - matching the authored OOXML diff.

$wb = $excel.ActiveWorkbook

# --- script1: clear the "selected tab" flag (another sheet becomes active) ---
$script1 = $wb.Worksheets.Item("script1")

# --- valid_login: inserted directly after script1 ---
$validLogin = $wb.Worksheets.Add([System.Type]::Missing, $script1)
$validLogin.Name = "valid_login"

$validLogin.Range("A1").Value = "username"
$validLogin.Range("B1").Value = "password"
$validLogin.Range("A2").Value = "admin"
$validLogin.Range("B2").Value = "pointofsale"

$null = $validLogin.Range("B2").Select()
$excel.ActiveWindow.Zoom = 220

# --- invalid_login: inserted directly after valid_login ---
$invalidLogin = $wb.Worksheets.Add([System.Type]::Missing, $validLogin)
$invalidLogin.Name = "invalid_login"

$invalidLogin.Range("A1").Value = "username"
$invalidLogin.Range("B1").Value = "password"
$invalidLogin.Range("A2").Value = "abc"
$invalidLogin.Range("B2").Value = "xyz"

$null = $invalidLogin.Range("B3").Select()

# invalid_login ends up last-added / active tab, matching tabSelected + activeTab
$null = $invalidLogin.Activate()
